# Applies the cryptos.xlsx crypto-price-table refresh described in the commit
# "Updated cryptos list on Fri Jan  5 19:25:37 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "43.852.54"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -1.34%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.236.95"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -2.35%  "
$ws.Range("E4").Value = "  +0.26%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "314.61"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.96%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "98.66"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -6.20%  "
$ws.Range("E7").Value = "  -3.25%  "
$ws.Range("E8").Value = "  +0.27%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.533"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -7.55%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "35.87"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -8.52%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0821"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -2.77%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "7.35"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -7.59%  "
$ws.Range("E13").Value = "  -2.96%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "2.576.68"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.838"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -5.45%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.237.93"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.91%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "13.95"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -5.09%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "43.704.22"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.33%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.09"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -8.49%  "
$ws.Range("E20").Value = "  -3.80%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.29"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -4.83%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "65.92"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.97%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "236.46"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("E24").Value = "  -7.74%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.03"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -8.57%  "
$ws.Range("E26").Value = "  +0.39%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.09"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("E28").Value = "  -3.19%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "36.56"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -6.93%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.98"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -9.22%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "19.98"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -3.18%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "156.48"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -4.43%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.0831"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -6.62%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.34"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("E36").Value = "  -9.16%  "
$ws.Range("E37").Value = "  -4.88%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.117"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -3.37%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "15.54"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  -11.06%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "3.99"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -12.19%  "
$ws.Range("E42").Value = "  -6.61%  "
$ws.Range("E43").Value = "  +0.33%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.705.01"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -4.00%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "82.48"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -4.32%  "
$ws.Range("E46").Value = "  -7.31%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "5.14"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -5.56%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "101.70"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "15.12"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "71.34"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -5.96%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "56.22"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -6.32%  "
